$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove G3 value
$ws.Range("G3").ClearContents()

# Remove F6 value
$ws.Range("F6").ClearContents()

# Add F7 value "fail"
$ws.Range("F7").Value = "fail"

# Update selection to F9
$ws.Range("F9").Select()
